$wb = $excel.ActiveWorkbook

# --- Settings sheet: change the TripAdvisorURL value cell (B2) to the new
#     GetYourGuide URL and turn it into a real hyperlink ---
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Hyperlinks.Add($wsSettings.Range("B2"), "https://www.getyourguide.com", "", "", "https://www.getyourguide.com")

# --- Make "Settings" the active/selected sheet (was "Assets" before) ---
$wsSettings.Activate() | Out-Null
$wsSettings.Range("B2").Select() | Out-Null
